$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.96"
$ws.Range("E2").Value = "'0.07%"
$ws.Range("D3").Value = "'38.41"
$ws.Range("E3").Value = "'7.58%"
$ws.Range("D4").Value = "'5.088"
$ws.Range("E4").Value = "'0.88%"
$ws.Range("D5").Value = "'0.08052"
$ws.Range("E5").Value = "'0.67%"
$ws.Range("D6").Value = "'1.929"
$ws.Range("E6").Value = "'3.41%"
$ws.Range("D7").Value = "'4.198"
$ws.Range("E7").Value = "'2.00%"
$ws.Range("D8").Value = "'7.945"
$ws.Range("E8").Value = "'2.22%"
$ws.Range("D9").Value = "'0.9300"
$ws.Range("E9").Value = "'0.95%"
$ws.Range("D10").Value = "'0.1433"
$ws.Range("E10").Value = "'10.99%"
$ws.Range("D11").Value = "'0.1924"
$ws.Range("D12").Value = "'0.09030"
$ws.Range("E12").Value = "'-0.73%"
$ws.Range("D13").Value = "'0.03526"
$ws.Range("E13").Value = "'3.03%"
$ws.Range("D14").Value = "'0.09769"
$ws.Range("E14").Value = "'-1.11%"
$ws.Range("D15").Value = "'0.001395"
$ws.Range("E15").Value = "'-1.40%"
$ws.Range("D16").Value = "'0.006051"
$ws.Range("E16").Value = "'-2.85%"
$ws.Range("D17").Value = "'3.725"
$ws.Range("E17").Value = "'-3.35%"
$ws.Range("E18").Value = "'0.95%"
$ws.Range("D19").Value = "'0.3461"
$ws.Range("E19").Value = "'1.29%"
$ws.Range("E20").Value = "'-0.41%"
$ws.Range("E21").Value = "'1.04%"
$ws.Range("E22").Value = "'-3.54%"
$ws.Range("D23").Value = "'0.04368"
$ws.Range("E23").Value = "'-1.03%"
$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'-0.22%"
$ws.Range("D25").Value = "'0.004115"
$ws.Range("E25").Value = "'-15.71%"
$ws.Range("E27").Value = "'0.23%"
$ws.Range("D39").Value = "'0.02082"
$ws.Range("E39").Value = "'7.25%"
$ws.Range("D40").Value = "'0.05023"
$ws.Range("E40").Value = "'-2.86%"
$ws.Range("D41").Value = "'0.007475"
$ws.Range("E41").Value = "'-1.18%"
$ws.Range("E42").Value = "'0.06%"
$ws.Range("D43").Value = "'0.1348"
$ws.Range("E43").Value = "'-0.39%"
$ws.Range("D44").Value = "'0.002146"
$ws.Range("E44").Value = "'0.70%"
$ws.Range("D45").Value = "'0.008901"
$ws.Range("E45").Value = "'-10.16%"
$ws.Range("D46").Value = "'0.00006189"
$ws.Range("E46").Value = "'0.03%"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("D48").Value = "'0.002986"
$ws.Range("D49").Value = "'0.001602"
$ws.Range("E49").Value = "'28.15%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.21%"
